$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''61.199.84'
$ws.Range("E2").Value = '''  +1.50%  '

# Row 3
$ws.Range("D3").Value = '''2.646.42'
$ws.Range("E3").Value = '''  +2.18%  '

# Row 4
$ws.Range("E4").Value = '''  -0.02%  '

# Row 5
$ws.Range("D5").Value = '''532.90'
$ws.Range("E5").Value = '''  +4.62%  '

# Row 6
$ws.Range("D6").Value = '''156.65'
$ws.Range("E6").Value = '''  +1.82%  '

# Row 7
$ws.Range("E7").Value = '''  -0.13%  '

# Row 8
$ws.Range("D8").Value = '''0.592'
$ws.Range("E8").Value = '''  +0.04%  '

# Row 9
$ws.Range("D9").Value = '''6.71'
$ws.Range("E9").Value = '''  +0.28%  '

# Row 10
$ws.Range("D10").Value = '''0.111'
$ws.Range("E10").Value = '''  +6.51%  '

# Row 11
$ws.Range("D11").Value = '''0.351'
$ws.Range("E11").Value = '''  +1.39%  '

# Row 12
$ws.Range("E12").Value = '''  +0.37%  '

# Row 13
$ws.Range("D13").Value = '''3.102.53'
$ws.Range("E13").Value = '''  +1.96%  '

# Row 14
$ws.Range("D14").Value = '''61.183.92'
$ws.Range("E14").Value = '''  +1.47%  '

# Row 15
$ws.Range("D15").Value = '''22.13'
$ws.Range("E15").Value = '''  +2.62%  '

# Row 16
$ws.Range("E16").Value = '''  +4.39%  '

# Row 17
$ws.Range("D17").Value = '''2.642.48'
$ws.Range("E17").Value = '''  +1.98%  '

# Row 18
$ws.Range("D18").Value = '''4.81'

# Row 19
$ws.Range("D19").Value = '''358.55'
$ws.Range("E19").Value = '''  +1.92%  '

# Row 20
$ws.Range("D20").Value = '''10.72'
$ws.Range("E20").Value = '''  +1.70%  '

# Row 21
$ws.Range("D21").Value = '''6.28'
$ws.Range("E21").Value = '''  +2.78%  '

# Row 22
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '''  +0.13%  '

# Row 23
$ws.Range("D23").Value = '''62.00'
$ws.Range("E23").Value = '''  +2.96%  '

# Row 24
$ws.Range("D24").Value = '''0.434'
$ws.Range("E24").Value = '''  +3.13%  '

# Row 25
$ws.Range("D25").Value = '''0.170'
$ws.Range("E25").Value = '''  +2.18%  '

# Row 26
$ws.Range("D26").Value = '''2.745.80'
$ws.Range("E26").Value = '''  +1.39%  '

# Row 27
$ws.Range("D27").Value = '''0.997'
$ws.Range("E27").Value = '''  -0.12%  '

# Row 28
$ws.Range("D28").Value = '''0.0₃0875'
$ws.Range("E28").Value = '''  +4.04%  '

# Row 29
$ws.Range("D29").Value = '''7.47'
$ws.Range("E29").Value = '''  +1.64%  '

# Row 30
$ws.Range("E30").Value = '''  -0.10%  '

# Row 31
$ws.Range("D31").Value = '''6.19'
$ws.Range("E31").Value = '''  +7.97%  '

# Row 32
$ws.Range("D32").Value = '''19.65'
$ws.Range("E32").Value = '''  +1.40%  '

# Row 33
$ws.Range("E33").Value = '''  +4.39%  '

# Row 34
$ws.Range("D34").Value = '''151.42'
$ws.Range("E34").Value = '''  -0.39%  '

# Row 35
$ws.Range("D35").Value = '''4.21'
$ws.Range("E35").Value = '''  +4.98%  '

# Row 36
$ws.Range("E36").Value = '''  +2.28%  '

# Row 37
$ws.Range("D37").Value = '''0.928'
$ws.Range("E37").Value = '''  +10.15%  '

# Row 38
$ws.Range("D38").Value = '''0.889'
$ws.Range("E38").Value = '''  +3.61%  '

# Row 39
$ws.Range("E39").Value = '''  +2.37%  '

# Row 40
$ws.Range("D40").Value = '''3.84'
$ws.Range("E40").Value = '''  +2.38%  '

# Row 41
$ws.Range("D41").Value = '''297.67'
$ws.Range("E41").Value = '''  -0.98%  '

# Row 42
$ws.Range("D42").Value = '''0.647'
$ws.Range("E42").Value = '''  +4.76%  '

# Row 43
$ws.Range("E43").Value = '''  +1.78%  '

# Row 44
$ws.Range("D44").Value = '''0.0565'
$ws.Range("E44").Value = '''  +2.20%  '

# Row 45
$ws.Range("D45").Value = '''5.14'
$ws.Range("E45").Value = '''  +6.94%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''20.04'
$ws.Range("E46").Value = '''  +1.99%  '

# Row 47
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").Value = '''0.997'
$ws.Range("E47").Value = '''  -0.02%  '

# Row 48
$ws.Range("D48").Value = '''0.0240'
$ws.Range("E48").Value = '''  +3.31%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '''19.22'
$ws.Range("E49").Value = '''  +7.26%  '

# Row 50
$ws.Range("D50").Value = '''10.36'
$ws.Range("E50").Value = '''  +0.54%  '

# Row 51
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '''1.87'
$ws.Range("E51").Value = '''  +5.49%  '

